$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D60").Value = 44526
$ws.Range("D61").Value = 44286
$ws.Range("J61").Value = 160
$ws.Range("K61").Value = 1500
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = 1500
$ws.Range("P61").Value = 500
$ws.Range("D62").Value = 44335
$ws.Range("J62").Value = 160
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 1500
$ws.Range("M62").Value = 1500
$ws.Range("P62").Value = 500
$ws.Range("D63").Value = 44186
$ws.Range("J63").Value = 180
$ws.Range("K63").Value = 1500
$ws.Range("L63").Value = 1500
$ws.Range("M63").Value = 1500
$ws.Range("P63").Value = 500
$ws.Range("D64").Value = 44460
$ws.Range("J64").Value = 160
$ws.Range("K64").Value = 1500
$ws.Range("L64").Value = 1500
$ws.Range("M64").Value = 1500
$ws.Range("P64").Value = 500
$ws.Range("D65").Value = 44438
$ws.Range("J65").Value = 160
$ws.Range("K65").Value = 1500
$ws.Range("L65").Value = 1500
$ws.Range("M65").Value = 1500
$ws.Range("P65").Value = 500
$ws.Range("D66").Value = 44519
$ws.Range("J66").Value = 160
$ws.Range("K66").Value = 1500
$ws.Range("L66").Value = 1500
$ws.Range("M66").Value = 1500
$ws.Range("P66").Value = 500
$ws.Range("D67").Value = 44392
$ws.Range("J67").Value = 160
$ws.Range("K67").Value = 1500
$ws.Range("L67").Value = 1500
$ws.Range("M67").Value = 1500
$ws.Range("P67").Value = 500
$ws.Range("D68").Value = 44355
$ws.Range("J68").Value = 180
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 1500
$ws.Range("M68").Value = 1500
$ws.Range("P68").Value = 500
$ws.Range("D69").Value = 44489
$ws.Range("J69").Value = 160
$ws.Range("K69").Value = 1500
$ws.Range("L69").Value = 1500
$ws.Range("M69").Value = 1500
$ws.Range("P69").Value = 500
$ws.Range("D70").Value = 44434
$ws.Range("J70").Value = 140
$ws.Range("K70").Value = 1500
$ws.Range("L70").Value = 1500
$ws.Range("M70").Value = 1500
$ws.Range("P70").Value = 500
$ws.Range("D71").Value = 44497
$ws.Range("J71").Value = 160
$ws.Range("K71").Value = 1500
$ws.Range("L71").Value = 1500
$ws.Range("M71").Value = 1500
$ws.Range("P71").Value = 500
$ws.Range("D72").Value = 44358
$ws.Range("J72").Value = 160
$ws.Range("K72").Value = 1500
$ws.Range("L72").Value = 1500
$ws.Range("M72").Value = 1500
$ws.Range("P72").Value = 500
$ws.Range("D73").Value = 44399
$ws.Range("J73").Value = 120
$ws.Range("K73").Value = 1500
$ws.Range("L73").Value = 1500
$ws.Range("M73").Value = 1500
$ws.Range("P73").Value = 500
$ws.Range("D74").Value = 44298
$ws.Range("J74").Value = 160
$ws.Range("K74").Value = 1500
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = 1500
$ws.Range("P74").Value = 500
$ws.Range("D75").Value = 44482
$ws.Range("J75").Value = 160
$ws.Range("K75").Value = 1500
$ws.Range("L75").Value = 1500
$ws.Range("M75").Value = 1500
$ws.Range("P75").Value = 500
$ws.Range("D76").Value = 44405
$ws.Range("J76").Value = 160
$ws.Range("K76").Value = 1500
$ws.Range("L76").Value = 1500
$ws.Range("M76").Value = 1500
$ws.Range("P76").Value = 500
$ws.Range("D77").Value = 44250
$ws.Range("J77").Value = 160
$ws.Range("K77").Value = 1500
$ws.Range("L77").Value = 1500
$ws.Range("M77").Value = 1500
$ws.Range("P77").Value = 500
$ws.Range("D78").Value = 44218
$ws.Range("J78").Value = 130
$ws.Range("K78").Value = 1500
$ws.Range("L78").Value = 1500
$ws.Range("M78").Value = 1500
$ws.Range("P78").Value = 500
$ws.Range("D79").Value = 44273
$ws.Range("J79").Value = 160
$ws.Range("K79").Value = 1500
$ws.Range("L79").Value = 1500
$ws.Range("M79").Value = 1500
$ws.Range("P79").Value = 500
$ws.Range("D80").Value = 44386
$ws.Range("J80").Value = 160
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 1500
$ws.Range("M80").Value = 1500
$ws.Range("P80").Value = 500
$ws.Range("D81").Value = 44435
$ws.Range("J81").Value = 810
$ws.Range("K81").Value = 1500
$ws.Range("L81").Value = 1500
$ws.Range("M81").Value = 1500
$ws.Range("P81").Value = 500
$ws.Range("D82").Value = 44328
$ws.Range("J82").Value = 160
$ws.Range("K82").Value = 1500
$ws.Range("L82").Value = 1500
$ws.Range("M82").Value = 1500
$ws.Range("P82").Value = 500
$ws.Range("D83").Value = 44277
$ws.Range("J83").Value = 160
$ws.Range("K83").Value = 1500
$ws.Range("L83").Value = 1500
$ws.Range("M83").Value = 1500
$ws.Range("P83").Value = 500
$ws.Range("D84").Value = 44442
$ws.Range("J84").Value = 180
$ws.Range("K84").Value = 1500
$ws.Range("L84").Value = 1500
$ws.Range("M84").Value = 1500
$ws.Range("P84").Value = 500
$ws.Range("D85").Value = 44516
$ws.Range("J85").Value = 150
$ws.Range("K85").Value = 1500
$ws.Range("L85").Value = 1500
$ws.Range("M85").Value = 1500
$ws.Range("P85").Value = 500
$ws.Range("D86").Value = 44175
$ws.Range("J86").Value = 120
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = 1500
$ws.Range("P86").Value = 500
$ws.Range("D87").Value = 44168
$ws.Range("J87").Value = 160
$ws.Range("K87").Value = 1500
$ws.Range("L87").Value = 1500
$ws.Range("M87").Value = 1500
$ws.Range("P87").Value = 500
$ws.Range("D88").Value = 44203
$ws.Range("J88").Value = 120
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 1500
$ws.Range("M88").Value = 1500
$ws.Range("P88").Value = 500
$ws.Range("D89").Value = 44475
$ws.Range("J89").Value = 160
$ws.Range("K89").Value = 1500
$ws.Range("L89").Value = 1500
$ws.Range("M89").Value = 1500
$ws.Range("P89").Value = 500
$ws.Range("D90").Value = 44483
$ws.Range("J90").Value = 180
$ws.Range("K90").Value = 1500
$ws.Range("L90").Value = 1500
$ws.Range("M90").Value = 1500
$ws.Range("P90").Value = 500
$ws.Range("D91").Value = 44217
$ws.Range("J91").Value = 120
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 1500
$ws.Range("M91").Value = 1500
$ws.Range("P91").Value = 500
$ws.Range("D92").Value = 44235
$ws.Range("J92").Value = 160
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = 1500
$ws.Range("P92").Value = 500
$ws.Range("D93").Value = 44200
$ws.Range("J93").Value = 120
$ws.Range("K93").Value = 1500
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = 1500
$ws.Range("P93").Value = 500
$ws.Range("D94").Value = 44419
$ws.Range("J94").Value = 130
$ws.Range("K94").Value = 1500
$ws.Range("L94").Value = 1500
$ws.Range("M94").Value = 1500
$ws.Range("P94").Value = 500
$ws.Range("D95").Value = 44162
$ws.Range("J95").Value = 160
$ws.Range("K95").Value = 1500
$ws.Range("L95").Value = 1500
$ws.Range("M95").Value = 1500
$ws.Range("P95").Value = 500
$ws.Range("D96").Value = 44357
$ws.Range("J96").Value = 160
$ws.Range("K96").Value = 1500
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = 1500
$ws.Range("P96").Value = 500
$ws.Range("D97").Value = 44244
$ws.Range("J97").Value = 110
$ws.Range("K97").Value = 1500
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = 1500
$ws.Range("P97").Value = 500
$ws.Range("D98").Value = 44202
$ws.Range("J98").Value = 120
$ws.Range("K98").Value = 1500
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = 1500
$ws.Range("P98").Value = 500
$ws.Range("D99").Value = 44333
$ws.Range("J99").Value = 120
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = 1500
$ws.Range("P99").Value = 500
$ws.Range("D100").Value = 44320
$ws.Range("J100").Value = 160
$ws.Range("K100").Value = 1500
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = 1500
$ws.Range("P100").Value = 500
$ws.Range("D101").Value = 44252
$ws.Range("J101").Value = 160
$ws.Range("K101").Value = 1500
$ws.Range("L101").Value = 1500
$ws.Range("M101").Value = 1500
$ws.Range("P101").Value = 500
$ws.Range("D102").Value = 44467
$ws.Range("J102").Value = 160
$ws.Range("K102").Value = 1500
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 1500
$ws.Range("P102").Value = 500
$ws.Range("D103").Value = 44264
$ws.Range("J103").Value = 120
$ws.Range("K103").Value = 1500
$ws.Range("L103").Value = 1500
$ws.Range("M103").Value = 1500
$ws.Range("P103").Value = 500
$ws.Range("D104").Value = 44214
$ws.Range("J104").Value = 110
$ws.Range("K104").Value = 1500
$ws.Range("L104").Value = 1500
$ws.Range("M104").Value = 1500
$ws.Range("P104").Value = 500
$ws.Range("D105").Value = 44167
$ws.Range("J105").Value = 150
$ws.Range("K105").Value = 1500
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 1500
$ws.Range("P105").Value = 500
$ws.Range("D106").Value = 44291
$ws.Range("J106").Value = 89
$ws.Range("K106").Value = 1800
$ws.Range("L106").Value = 1800
$ws.Range("M106").Value = 1800
$ws.Range("P106").Value = 600
$ws.Range("D107").Value = 44174
$ws.Range("J107").Value = 180
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 1500
$ws.Range("P107").Value = 500
$ws.Range("D108").Value = 44293
$ws.Range("J108").Value = 160
$ws.Range("K108").Value = 1500
$ws.Range("L108").Value = 1500
$ws.Range("M108").Value = 1500
$ws.Range("P108").Value = 500
$ws.Range("D109").Value = 44496
$ws.Range("J109").Value = 150
$ws.Range("K109").Value = 1500
$ws.Range("L109").Value = 1500
$ws.Range("M109").Value = 1500
$ws.Range("P109").Value = 500
$ws.Range("D110").Value = 44326
$ws.Range("J110").Value = 120
$ws.Range("K110").Value = 1500
$ws.Range("L110").Value = 1500
$ws.Range("M110").Value = 1500
$ws.Range("P110").Value = 500
$ws.Range("D111").Value = 44302
$ws.Range("J111").Value = 130
$ws.Range("K111").Value = 1500
$ws.Range("L111").Value = 1500
$ws.Range("M111").Value = 1500
$ws.Range("P111").Value = 500
$ws.Range("D112").Value = 44308
$ws.Range("J112").Value = 160
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 1500
$ws.Range("M112").Value = 1500
$ws.Range("P112").Value = 500
$ws.Range("D113").Value = 44498
$ws.Range("J113").Value = 160
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1500
$ws.Range("P113").Value = 500
$ws.Range("D114").Value = 44420
$ws.Range("J114").Value = 160
$ws.Range("K114").Value = 1500
$ws.Range("L114").Value = 1500
$ws.Range("M114").Value = 1500
$ws.Range("P114").Value = 500
$ws.Range("D115").Value = 44398
$ws.Range("J115").Value = 160
$ws.Range("K115").Value = 1500
$ws.Range("L115").Value = 1500
$ws.Range("M115").Value = 1500
$ws.Range("P115").Value = 500
$ws.Range("D116").Value = 44396
$ws.Range("J116").Value = 160
$ws.Range("K116").Value = 1500
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = 1500
$ws.Range("P116").Value = 500
$ws.Range("D117").Value = 44321
$ws.Range("J117").Value = 130
$ws.Range("K117").Value = 1500
$ws.Range("L117").Value = 1500
$ws.Range("M117").Value = 1500
$ws.Range("P117").Value = 500
$ws.Range("D118").Value = 44208
$ws.Range("J118").Value = 160
$ws.Range("K118").Value = 1500
$ws.Range("L118").Value = 1500
$ws.Range("M118").Value = 1500
$ws.Range("P118").Value = 500
$ws.Range("D119").Value = 44349
$ws.Range("J119").Value = 160
$ws.Range("K119").Value = 1500
$ws.Range("L119").Value = 1500
$ws.Range("M119").Value = 1500
$ws.Range("P119").Value = 500
$ws.Range("D120").Value = 44477
$ws.Range("J120").Value = 160
$ws.Range("K120").Value = 1500
$ws.Range("L120").Value = 1500
$ws.Range("M120").Value = 1500
$ws.Range("P120").Value = 500
$ws.Range("D121").Value = 44487
$ws.Range("J121").Value = 160
$ws.Range("K121").Value = 1500
$ws.Range("L121").Value = 1500
$ws.Range("M121").Value = 1500
$ws.Range("P121").Value = 500
$ws.Range("D122").Value = 44452
$ws.Range("J122").Value = 190
$ws.Range("K122").Value = 1500
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = 1500
$ws.Range("P122").Value = 500
$ws.Range("D123").Value = 44211
$ws.Range("J123").Value = 120
$ws.Range("K123").Value = 1500
$ws.Range("L123").Value = 1500
$ws.Range("M123").Value = 1500
$ws.Range("P123").Value = 500
$ws.Range("D124").Value = 44505
$ws.Range("J124").Value = 120
$ws.Range("K124").Value = 1500
$ws.Range("L124").Value = 1500
$ws.Range("M124").Value = 1500
$ws.Range("P124").Value = 500
$ws.Range("D125").Value = 44204
$ws.Range("J125").Value = 180
$ws.Range("K125").Value = 1500
$ws.Range("L125").Value = 1500
$ws.Range("M125").Value = 1500
$ws.Range("P125").Value = 500
$ws.Range("D126").Value = 44306
$ws.Range("J126").Value = 160
$ws.Range("K126").Value = 1500
$ws.Range("L126").Value = 1500
$ws.Range("M126").Value = 1500
$ws.Range("P126").Value = 500
$ws.Range("D127").Value = 44509
$ws.Range("J127").Value = 160
$ws.Range("K127").Value = 1500
$ws.Range("L127").Value = 1500
$ws.Range("M127").Value = 1500
$ws.Range("P127").Value = 500
$ws.Range("D128").Value = 44454
$ws.Range("J128").Value = 160
$ws.Range("K128").Value = 1500
$ws.Range("L128").Value = 1500
$ws.Range("M128").Value = 1500
$ws.Range("P128").Value = 500
$ws.Range("D129").Value = 44189
$ws.Range("J129").Value = 180
$ws.Range("K129").Value = 1500
$ws.Range("L129").Value = 1500
$ws.Range("M129").Value = 1500
$ws.Range("P129").Value = 500
$ws.Range("D130").Value = 44278
$ws.Range("J130").Value = 130
$ws.Range("K130").Value = 1500
$ws.Range("L130").Value = 1500
$ws.Range("M130").Value = 1500
$ws.Range("P130").Value = 500
$ws.Range("D131").Value = 44265
$ws.Range("J131").Value = 120
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 1500
$ws.Range("M131").Value = 1500
$ws.Range("P131").Value = 500
$ws.Range("D132").Value = 44494
$ws.Range("J132").Value = 190
$ws.Range("K132").Value = 1500
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = 1500
$ws.Range("P132").Value = 500
$ws.Range("D133").Value = 44300
$ws.Range("J133").Value = 160
$ws.Range("K133").Value = 1500
$ws.Range("L133").Value = 1500
$ws.Range("M133").Value = 1500
$ws.Range("P133").Value = 500
$ws.Range("D134").Value = 44209
$ws.Range("J134").Value = 160
$ws.Range("K134").Value = 1500
$ws.Range("L134").Value = 1500
$ws.Range("M134").Value = 1500
$ws.Range("P134").Value = 500
$ws.Range("D135").Value = 44237
$ws.Range("J135").Value = 130
$ws.Range("K135").Value = 1500
$ws.Range("L135").Value = 1500
$ws.Range("M135").Value = 1500
$ws.Range("P135").Value = 500
$ws.Range("D136").Value = 44356
$ws.Range("J136").Value = 160
$ws.Range("K136").Value = 1500
$ws.Range("L136").Value = 1500
$ws.Range("M136").Value = 1500
$ws.Range("P136").Value = 500
$ws.Range("D137").Value = 44469
$ws.Range("J137").Value = 160
$ws.Range("K137").Value = 1500
$ws.Range("L137").Value = 1500
$ws.Range("M137").Value = 1500
$ws.Range("P137").Value = 500
$ws.Range("D138").Value = 44453
$ws.Range("J138").Value = 130
$ws.Range("K138").Value = 1500
$ws.Range("L138").Value = 1500
$ws.Range("M138").Value = 1500
$ws.Range("P138").Value = 500
$ws.Range("D139").Value = 44518
$ws.Range("J139").Value = 160
$ws.Range("K139").Value = 1500
$ws.Range("L139").Value = 1500
$ws.Range("M139").Value = 1500
$ws.Range("P139").Value = 500
$ws.Range("D140").Value = 44446
$ws.Range("J140").Value = 180
$ws.Range("K140").Value = 1500
$ws.Range("L140").Value = 1500
$ws.Range("M140").Value = 1500
$ws.Range("P140").Value = 500
$ws.Range("D141").Value = 44463
$ws.Range("J141").Value = 160
$ws.Range("K141").Value = 1500
$ws.Range("L141").Value = 1500
$ws.Range("M141").Value = 1500
$ws.Range("P141").Value = 500
$ws.Range("D142").Value = 44245
$ws.Range("J142").Value = 120
$ws.Range("K142").Value = 1500
$ws.Range("L142").Value = 1500
$ws.Range("M142").Value = 1500
$ws.Range("P142").Value = 500
$ws.Range("D143").Value = 44323
$ws.Range("J143").Value = 160
$ws.Range("K143").Value = 1500
$ws.Range("L143").Value = 1500
$ws.Range("M143").Value = 1500
$ws.Range("P143").Value = 500
$ws.Range("D144").Value = 44229
$ws.Range("J144").Value = 160
$ws.Range("K144").Value = 1500
$ws.Range("L144").Value = 1500
$ws.Range("M144").Value = 1500
$ws.Range("P144").Value = 500
$ws.Range("D145").Value = 44417
$ws.Range("J145").Value = 160
$ws.Range("K145").Value = 1500
$ws.Range("L145").Value = 1500
$ws.Range("M145").Value = 1500
$ws.Range("P145").Value = 500
$ws.Range("D146").Value = 44445
$ws.Range("J146").Value = 180
$ws.Range("K146").Value = 1500
$ws.Range("L146").Value = 1500
$ws.Range("M146").Value = 1500
$ws.Range("P146").Value = 500
$ws.Range("D147").Value = 44249
$ws.Range("J147").Value = 160
$ws.Range("K147").Value = 1500
$ws.Range("L147").Value = 1500
$ws.Range("M147").Value = 1500
$ws.Range("P147").Value = 500
$ws.Range("D148").Value = 44342
$ws.Range("J148").Value = 260
$ws.Range("K148").Value = 1500
$ws.Range("L148").Value = 1500
$ws.Range("M148").Value = 1500
$ws.Range("P148").Value = 500
$ws.Range("D149").Value = 44523
$ws.Range("J149").Value = 160
$ws.Range("K149").Value = 1500
$ws.Range("L149").Value = 1500
$ws.Range("M149").Value = 1500
$ws.Range("P149").Value = 500
$ws.Range("D150").Value = 44259
$ws.Range("J150").Value = 120
$ws.Range("K150").Value = 1500
$ws.Range("L150").Value = 1500
$ws.Range("M150").Value = 1500
$ws.Range("P150").Value = 500
$ws.Range("D151").Value = 44216
$ws.Range("J151").Value = 80
$ws.Range("K151").Value = 1500
$ws.Range("L151").Value = 1500
$ws.Range("M151").Value = 1500
$ws.Range("P151").Value = 500
$ws.Range("D152").Value = 44406
$ws.Range("J152").Value = 160
$ws.Range("K152").Value = 1500
$ws.Range("L152").Value = 1500
$ws.Range("M152").Value = 1500
$ws.Range("P152").Value = 500
$ws.Range("D153").Value = 44295
$ws.Range("J153").Value = 120
$ws.Range("K153").Value = 1500
$ws.Range("L153").Value = 1500
$ws.Range("M153").Value = 1500
$ws.Range("P153").Value = 500
$ws.Range("D154").Value = 44270
$ws.Range("J154").Value = 120
$ws.Range("K154").Value = 1500
$ws.Range("L154").Value = 1500
$ws.Range("M154").Value = 1500
$ws.Range("P154").Value = 500
$ws.Range("D155").Value = 44363
$ws.Range("J155").Value = 130
$ws.Range("K155").Value = 1500
$ws.Range("L155").Value = 1500
$ws.Range("M155").Value = 1500
$ws.Range("P155").Value = 500
$ws.Range("D156").Value = 44299
$ws.Range("J156").Value = 130
$ws.Range("K156").Value = 1500
$ws.Range("L156").Value = 1500
$ws.Range("M156").Value = 1500
$ws.Range("P156").Value = 500
$ws.Range("D157").Value = 44257
$ws.Range("J157").Value = 120
$ws.Range("K157").Value = 1500
$ws.Range("L157").Value = 1500
$ws.Range("M157").Value = 1500
$ws.Range("P157").Value = 500
$ws.Range("D158").Value = 44336
$ws.Range("J158").Value = 160
$ws.Range("K158").Value = 1500
$ws.Range("L158").Value = 1500
$ws.Range("M158").Value = 1500
$ws.Range("P158").Value = 500
$ws.Range("D159").Value = 44372
$ws.Range("J159").Value = 160
$ws.Range("K159").Value = 1500
$ws.Range("L159").Value = 1500
$ws.Range("M159").Value = 1500
$ws.Range("P159").Value = 500
$ws.Range("D160").Value = 44403
$ws.Range("J160").Value = 180
$ws.Range("K160").Value = 1500
$ws.Range("L160").Value = 1500
$ws.Range("M160").Value = 1500
$ws.Range("P160").Value = 500
$ws.Range("D161").Value = 44195
$ws.Range("J161").Value = 180
$ws.Range("K161").Value = 1500
$ws.Range("L161").Value = 1500
$ws.Range("M161").Value = 1500
$ws.Range("P161").Value = 500
$ws.Range("D162").Value = 44376
$ws.Range("J162").Value = 160
$ws.Range("K162").Value = 1500
$ws.Range("L162").Value = 1500
$ws.Range("M162").Value = 1500
$ws.Range("P162").Value = 500
$ws.Range("D163").Value = 44474
$ws.Range("J163").Value = 160
$ws.Range("K163").Value = 1500
$ws.Range("L163").Value = 1500
$ws.Range("M163").Value = 1500
$ws.Range("P163").Value = 500
$ws.Range("D164").Value = 44524
$ws.Range("J164").Value = 160
$ws.Range("K164").Value = 1500
$ws.Range("L164").Value = 1500
$ws.Range("M164").Value = 1500
$ws.Range("P164").Value = 500
$ws.Range("D165").Value = 44172
$ws.Range("J165").Value = 110
$ws.Range("K165").Value = 1500
$ws.Range("L165").Value = 1500
$ws.Range("M165").Value = 1500
$ws.Range("P165").Value = 500
$ws.Range("D166").Value = 44421
$ws.Range("J166").Value = 180
$ws.Range("K166").Value = 1500
$ws.Range("L166").Value = 1500
$ws.Range("M166").Value = 1500
$ws.Range("P166").Value = 500
$ws.Range("D167").Value = 44431
$ws.Range("J167").Value = 180
$ws.Range("K167").Value = 1500
$ws.Range("L167").Value = 1500
$ws.Range("M167").Value = 1500
$ws.Range("P167").Value = 500
$ws.Range("D168").Value = 44239
$ws.Range("J168").Value = 120
$ws.Range("K168").Value = 1500
$ws.Range("L168").Value = 1500
$ws.Range("M168").Value = 1500
$ws.Range("P168").Value = 500
$ws.Range("D169").Value = 44426
$ws.Range("J169").Value = 160
$ws.Range("K169").Value = 1500
$ws.Range("L169").Value = 1500
$ws.Range("M169").Value = 1500
$ws.Range("P169").Value = 500
$ws.Range("D170").Value = 44448
$ws.Range("J170").Value = 160
$ws.Range("K170").Value = 1500
$ws.Range("L170").Value = 1500
$ws.Range("M170").Value = 1500
$ws.Range("P170").Value = 500
$ws.Range("D171").Value = 44362
$ws.Range("J171").Value = 180
$ws.Range("K171").Value = 1500
$ws.Range("L171").Value = 1500
$ws.Range("M171").Value = 1500
$ws.Range("P171").Value = 500
$ws.Range("D172").Value = 44210
$ws.Range("J172").Value = 120
$ws.Range("K172").Value = 1500
$ws.Range("L172").Value = 1500
$ws.Range("M172").Value = 1500
$ws.Range("P172").Value = 500
$ws.Range("D173").Value = 44176
$ws.Range("J173").Value = 80
$ws.Range("K173").Value = 1500
$ws.Range("L173").Value = 1500
$ws.Range("M173").Value = 1500
$ws.Range("P173").Value = 500
$ws.Range("D174").Value = 44301
$ws.Range("J174").Value = 130
$ws.Range("K174").Value = 1500
$ws.Range("L174").Value = 1500
$ws.Range("M174").Value = 1500
$ws.Range("P174").Value = 500
$ws.Range("D175").Value = 44407
$ws.Range("J175").Value = 160
$ws.Range("K175").Value = 1500
$ws.Range("L175").Value = 1500
$ws.Range("M175").Value = 1500
$ws.Range("P175").Value = 500
$ws.Range("D176").Value = 44284
$ws.Range("J176").Value = 180
$ws.Range("K176").Value = 1500
$ws.Range("L176").Value = 1500
$ws.Range("M176").Value = 1500
$ws.Range("P176").Value = 500
$ws.Range("D177").Value = 44441
$ws.Range("J177").Value = 190
$ws.Range("K177").Value = 1500
$ws.Range("L177").Value = 1500
$ws.Range("M177").Value = 1500
$ws.Range("P177").Value = 500
$ws.Range("D178").Value = 44279
$ws.Range("J178").Value = 160
$ws.Range("K178").Value = 1500
$ws.Range("L178").Value = 1500
$ws.Range("M178").Value = 1500
$ws.Range("P178").Value = 500
$ws.Range("D179").Value = 44341
$ws.Range("J179").Value = 160
$ws.Range("K179").Value = 1500
$ws.Range("L179").Value = 1500
$ws.Range("M179").Value = 1500
$ws.Range("P179").Value = 500
$ws.Range("D180").Value = 44504
$ws.Range("J180").Value = 160
$ws.Range("K180").Value = 1500
$ws.Range("L180").Value = 1500
$ws.Range("M180").Value = 1500
$ws.Range("P180").Value = 500
$ws.Range("D181").Value = 44350
$ws.Range("J181").Value = 160
$ws.Range("K181").Value = 1500
$ws.Range("L181").Value = 1500
$ws.Range("M181").Value = 1500
$ws.Range("P181").Value = 500
$ws.Range("D182").Value = 44312
$ws.Range("J182").Value = 160
$ws.Range("K182").Value = 1500
$ws.Range("L182").Value = 1500
$ws.Range("M182").Value = 1500
$ws.Range("P182").Value = 500
$ws.Range("D183").Value = 44382
$ws.Range("J183").Value = 160
$ws.Range("K183").Value = 1500
$ws.Range("L183").Value = 1500
$ws.Range("M183").Value = 1500
$ws.Range("P183").Value = 500
$ws.Range("D184").Value = 44384
$ws.Range("J184").Value = 160
$ws.Range("K184").Value = 1500
$ws.Range("L184").Value = 1500
$ws.Range("M184").Value = 1500
$ws.Range("P184").Value = 500
$ws.Range("D185").Value = 44329
$ws.Range("J185").Value = 160
$ws.Range("K185").Value = 1500
$ws.Range("L185").Value = 1500
$ws.Range("M185").Value = 1500
$ws.Range("P185").Value = 500
$ws.Range("D186").Value = 44522
$ws.Range("J186").Value = 160
$ws.Range("K186").Value = 1500
$ws.Range("L186").Value = 1500
$ws.Range("M186").Value = 1500
$ws.Range("P186").Value = 500
$ws.Range("D187").Value = 44246
$ws.Range("J187").Value = 160
$ws.Range("K187").Value = 1500
$ws.Range("L187").Value = 1500
$ws.Range("M187").Value = 1500
$ws.Range("P187").Value = 500
$ws.Range("D188").Value = 44491
$ws.Range("J188").Value = 160
$ws.Range("K188").Value = 1500
$ws.Range("L188").Value = 1500
$ws.Range("M188").Value = 1500
$ws.Range("P188").Value = 500
$ws.Range("D189").Value = 44272
$ws.Range("J189").Value = 160
$ws.Range("K189").Value = 1500
$ws.Range("L189").Value = 1500
$ws.Range("M189").Value = 1500
$ws.Range("P189").Value = 500
$ws.Range("D190").Value = 44305
$ws.Range("J190").Value = 180
$ws.Range("K190").Value = 1500
$ws.Range("L190").Value = 1500
$ws.Range("M190").Value = 1500
$ws.Range("P190").Value = 500
$ws.Range("D191").Value = 44447
$ws.Range("J191").Value = 160
$ws.Range("K191").Value = 1500
$ws.Range("L191").Value = 1500
$ws.Range("M191").Value = 1500
$ws.Range("P191").Value = 500
$ws.Range("D192").Value = 44425
$ws.Range("J192").Value = 160
$ws.Range("K192").Value = 1500
$ws.Range("L192").Value = 1500
$ws.Range("M192").Value = 1500
$ws.Range("P192").Value = 500
$ws.Range("D193").Value = 44315
$ws.Range("J193").Value = 130
$ws.Range("K193").Value = 1500
$ws.Range("L193").Value = 1500
$ws.Range("M193").Value = 1500
$ws.Range("P193").Value = 500
$ws.Range("D194").Value = 44348
$ws.Range("J194").Value = 160
$ws.Range("K194").Value = 1500
$ws.Range("L194").Value = 1500
$ws.Range("M194").Value = 1500
$ws.Range("P194").Value = 500
$ws.Range("D195").Value = 44322
$ws.Range("J195").Value = 130
$ws.Range("K195").Value = 1500
$ws.Range("L195").Value = 1500
$ws.Range("M195").Value = 1500
$ws.Range("P195").Value = 500
$ws.Range("D196").Value = 44495
$ws.Range("J196").Value = 160
$ws.Range("K196").Value = 1500
$ws.Range("L196").Value = 1500
$ws.Range("M196").Value = 1500
$ws.Range("P196").Value = 500
$ws.Range("D197").Value = 44232
$ws.Range("J197").Value = 120
$ws.Range("K197").Value = 1500
$ws.Range("L197").Value = 1500
$ws.Range("M197").Value = 1500
$ws.Range("P197").Value = 500
$ws.Range("D198").Value = 44327
$ws.Range("J198").Value = 190
$ws.Range("K198").Value = 1500
$ws.Range("L198").Value = 1500
$ws.Range("M198").Value = 1500
$ws.Range("P198").Value = 500
$ws.Range("D199").Value = 44510
$ws.Range("J199").Value = 160
$ws.Range("K199").Value = 1500
$ws.Range("L199").Value = 1500
$ws.Range("M199").Value = 1500
$ws.Range("P199").Value = 500
$ws.Range("D200").Value = 44161
$ws.Range("J200").Value = 180
$ws.Range("K200").Value = 1500
$ws.Range("L200").Value = 1500
$ws.Range("M200").Value = 1500
$ws.Range("P200").Value = 500
$ws.Range("D201").Value = 44468
$ws.Range("J201").Value = 180
$ws.Range("K201").Value = 1500
$ws.Range("L201").Value = 1500
$ws.Range("M201").Value = 1500
$ws.Range("P201").Value = 500
$ws.Range("D202").Value = 44517
$ws.Range("J202").Value = 160
$ws.Range("K202").Value = 1500
$ws.Range("L202").Value = 1500
$ws.Range("M202").Value = 1500
$ws.Range("P202").Value = 500
$ws.Range("D203").Value = 44238
$ws.Range("J203").Value = 130
$ws.Range("K203").Value = 1500
$ws.Range("L203").Value = 1500
$ws.Range("M203").Value = 1500
$ws.Range("P203").Value = 500
$ws.Range("D204").Value = 44391
$ws.Range("J204").Value = 160
$ws.Range("K204").Value = 1500
$ws.Range("L204").Value = 1500
$ws.Range("M204").Value = 1500
$ws.Range("P204").Value = 500
$ws.Range("D205").Value = 44236
$ws.Range("J205").Value = 120
$ws.Range("K205").Value = 1500
$ws.Range("L205").Value = 1500
$ws.Range("M205").Value = 1500
$ws.Range("P205").Value = 500
$ws.Range("D206").Value = 44251
$ws.Range("J206").Value = 80
$ws.Range("K206").Value = 1500
$ws.Range("L206").Value = 1500
$ws.Range("M206").Value = 1500
$ws.Range("P206").Value = 500
$ws.Range("D207").Value = 44515
$ws.Range("J207").Value = 160
$ws.Range("K207").Value = 1500
$ws.Range("L207").Value = 1500
$ws.Range("M207").Value = 1500
$ws.Range("P207").Value = 500
$ws.Range("D208").Value = 44330
$ws.Range("J208").Value = 160
$ws.Range("K208").Value = 1500
$ws.Range("L208").Value = 1500
$ws.Range("M208").Value = 1500
$ws.Range("P208").Value = 500
$ws.Range("D209").Value = 44432
$ws.Range("J209").Value = 150
$ws.Range("K209").Value = 1500
$ws.Range("L209").Value = 1500
$ws.Range("M209").Value = 1500
$ws.Range("P209").Value = 500
$ws.Range("D210").Value = 44181
$ws.Range("J210").Value = 90
$ws.Range("K210").Value = 1500
$ws.Range("L210").Value = 1500
$ws.Range("M210").Value = 1500
$ws.Range("P210").Value = 500
$ws.Range("D211").Value = 44194
$ws.Range("J211").Value = 80
$ws.Range("K211").Value = 1500
$ws.Range("L211").Value = 1500
$ws.Range("M211").Value = 1500
$ws.Range("P211").Value = 500
$ws.Range("D212").Value = 44271
$ws.Range("J212").Value = 180
$ws.Range("K212").Value = 1500
$ws.Range("L212").Value = 1500
$ws.Range("M212").Value = 1500
$ws.Range("P212").Value = 500
$ws.Range("D213").Value = 44307
$ws.Range("J213").Value = 130
$ws.Range("K213").Value = 1500
$ws.Range("L213").Value = 1500
$ws.Range("M213").Value = 1500
$ws.Range("P213").Value = 500
$ws.Range("D214").Value = 44400
$ws.Range("J214").Value = 160
$ws.Range("K214").Value = 1500
$ws.Range("L214").Value = 1500
$ws.Range("M214").Value = 1500
$ws.Range("P214").Value = 500
$ws.Range("D215").Value = 44309
$ws.Range("J215").Value = 160
$ws.Range("K215").Value = 1500
$ws.Range("L215").Value = 1500
$ws.Range("M215").Value = 1500
$ws.Range("P215").Value = 500
$ws.Range("D216").Value = 44508
$ws.Range("J216").Value = 160
$ws.Range("K216").Value = 1500
$ws.Range("L216").Value = 1500
$ws.Range("M216").Value = 1500
$ws.Range("P216").Value = 500
$ws.Range("A217").Value = 3
$ws.Range("B217").Value = 'Femacal de La Calera'
$ws.Range("C217").Value = 'Coquimbo'
$ws.Range("D217").Value = 44201
$ws.Range("E217").Value = 5
$ws.Range("F217").Value = 100112039
$ws.Range("G217").Value = 'Ciboulette'
$ws.Range("H217").Value = 'Sin especificar'
$ws.Range("I217").Value = 'Primera'
$ws.Range("J217").Value = 120
$ws.Range("K217").Value = 1500
$ws.Range("L217").Value = 1500
$ws.Range("M217").Value = 1500
$ws.Range("N217").Value = '$/docena de atados'
$ws.Range("O217").Value = 'Provincia de Quillota'
$ws.Range("P217").Value = 500
$ws.Range("Q217").Value = 3
$ws.Range("R217").Value = 'Hortaliza'
$ws.Range("D217").NumberFormat = "YYYY-MM-DD HH:MM:SS"
